$d = $word.ActiveDocument

# 1. Title on the cover page
$d.Content.Find.Execute("Release Notes v4.0.5.4 Fix 76", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Release Notes v4.2.0.0 Fix 56", 2)

# 2. Introduction paragraph
$d.Content.Find.Execute("This document defines the changes made to the Network Manager product for fix release v4.0.5.4 Fix 76 and is specifically targeted at end users.  ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "This document defines the changes made to the Network Manager product for fix release v4.2.0.0 Fix 56 and is specifically targeted at end users.  ", 2)

# 3. Fix Description cell - first part (Network Manager ... Patchset.) - do this
#    before the lone "4.0.5.4" replacement below so the longer phrase can still
#    be matched in full.
$d.Content.Find.Execute("Network Manager 4.0.5.4 Fix 76 Patchset.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Network Manager 4.2.0.0 Fix 56 Patchset.", 2)

# 4. Fix Description cell - second part (Set flags ... batch interface loader)
$d.Content.Find.Execute("Set flags to ensure the same behaviour when loading CSV files as when using the batch interface loader ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Nm3gaz_qry changes to fix problems in PBI Query where no network found ", 2)

# 5. Baseline Release value in table
$d.Content.Find.Execute("4.0.5.4", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "4.2.0.0", 2)

# 6. Unzip instructions
$d.Content.Find.Execute("Unzip nm_4054_fix76.zip to a staging folder.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Unzip nm_4200_fix56.zip to a staging folder.", 2)

# 7. START sql prompt instructions (avoid touching the surrounding quote
#    characters so Word's smart-quote AutoFormat doesn't curl them)
$d.Content.Find.Execute("START nm_4054_fix76.sql", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "START nm_4200_fix56.sql", 2)

# 8. log output sql filename
$d.Content.Find.Execute("log_nm_4054_fix76.sql", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "log_nm_4200_fix56.sql", 2)

# 9. package filename
$d.Content.Find.Execute("nm3mapcapture_ins_inv.pkw", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "nm3gaz_qry.pkw", 2)

# 10. package version number
$d.Content.Find.Execute("2.5.1.1", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2.3.1.2", 2)

# 11. installer sql filename
$d.Content.Find.Execute("nm_4054_fix76.sql", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "nm_4200_fix56.sql", 2)

# 12. "Log No. Summary" chapter intro paragraph (merges split runs, text unchanged)
$d.Content.Find.Execute("This chapter summarises all software changes that have been made in this release. ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "This chapter summarises all software changes that have been made in this release. ", 2)

# 13. Log number
$d.Content.Find.Execute("0111859", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "0111785", 2)

# 14. Log description text
$d.Content.Find.Execute("Using this fix allows the loading of hierarchies of assets to complete when using the CSV loader. Flags that were currently set only during the batch loading are now configured inside the nm3mapcapture_ins_inv.ins_inv procedure so that subordinate metadata is properly handled.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "This fix repairs the problem with using PBI queries which executed to give no results with an error suggesting that no network obeying the criteria could be found", 2)

# 15. Tracking id
$d.Content.Find.Execute("8001314953", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "8001274633", 2)
